# Generate Report for Handoff
#
# Refreshes the handoff-generation timestamps for the
# 93c84428-ec77-4f42-8d59-1aab91cfa317.md row (row 7) across the
# Overview / zh-cn / de-de sheets, as produced by a fresh handoff run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-19 20:47:00"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-19 20:46:56"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-19 20:47:00"
